$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.019.29'
$ws.Range('E2').Value = '  +4.02%  '
$ws.Range('D3').Value = '2.251.42'
$ws.Range('E3').Value = '  +3.20%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.614'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.10'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.85%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('E9').Value = '  +6.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0937'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.97'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.101'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').Value = '2.584.83'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.63'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.50%  '
$ws.Range('D16').Value = '2.240.92'
$ws.Range('E16').Value = '  +2.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.808'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '42.929.51'
$ws.Range('E18').Value = '  +4.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000106'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '231.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.91%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +14.00%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.59%  '
$ws.Range('E28').Value = '  +2.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '38.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +27.41%  '
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.66'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.34'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0797'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('E34').Value = '  +4.27%  '
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.110'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0333'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +17.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.52%  '
$ws.Range('E40').Value = '  +3.01%  '
$ws.Range('E41').Value = '  +2.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.203'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '60.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '105.87'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0995'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.459'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +24.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.67%  '
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('E50').Value = '  +1.67%  '
$ws.Range('D51').Value = '2.458.40'
$ws.Range('E51').Value = '  +3.06%  '
